$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Patient name (row 6) ---
$ws.Range("A6").Value = "REYES"
$ws.Range("C6").Value = "BOLÒS"
$ws.Range("E6").Value = "ANABELLA"
$ws.Range("G6").Value = ""
$ws.Range("I6").Value = "2009---122761/201761928"

# --- Dirección actual (row 8) ---
$ws.Range("A8").Value = "6TA CALLE 5-20 "
$ws.Range("D8").Value = "SAN JOSE PINULA Z, 2"
$ws.Range("F8").Value = "SAN JOSE PINULA"
# H8 (Departamento = GUATEMALA) unchanged
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = "51121930"

# --- Fecha de nacimiento / Edad / Lugar de nacimiento / Sexo (row 12) ---
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "1982-07-31"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "35"
$ws.Range("H12").Value = "SAN JOSE PINULA"
# J12 (Sexo = Femenino) unchanged

# --- Estado Civil / Ocupación / Nacionalidad / No. de Cédula (row 14) ---
$ws.Range("A14").Value = "Casado"
# D14 (Ocupación = AMA DE CASA) unchanged, F14 (Nacionalidad = GUATEMALTECA) unchanged
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "2206688410103"

# --- Nombre del Padre / Nombre de la Madre (row 18) ---
$ws.Range("A18").Value = "JUAN DE JESUS REYES"
$ws.Range("F18").Value = "ENRIQUETA BOLÒS"

# --- Emergencia: nombre / parentesco (row 20) ---
$ws.Range("A20").Value = "mario rene davila"
$ws.Range("F20").Value = "esposo"

# --- Hora de ingreso / Servicio (row 24) ---
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "11:49:20"
$ws.Range("D24").Value = "LYP"
